{"js": "// Update the job-title/date lines for the three most recent positions and\n// the in-progress Master's degree line in the \"Exp\u00e9rience professionnelle\"\n// and \"Education\" sections of the r\u00e9sum\u00e9.\nconst replacements = [\n  {\n    oldText: \"Concepteur d\\u2019animation (janvier 2021 \\u00e0 ce jour)\",\n    newText: \"Animation Spark : Concepteur d\\u2019animations (Jan 2021 - Pr\\u00e9sent)\"\n  },\n  {\n    oldText: \"Concepteur d\\u2019animation (janvier 2018 \\u00e0 d\\u00e9cembre 2020)\",\n    newText: \"Pixel Studio : Concepteur d\\u2019animations (juin 2018 - D\\u00e9c 2020)\"\n  },\n  {\n    oldText: \"Concepteur d\\u2019animation junior (septembre 2016 \\u00e0 mai 2018)\",\n    newText: \"Animation flash : Concepteur d\\u2019animation junior (sep 2016 - mai 2018)\"\n  },\n  {\n    oldText: \"Master en arts d\\u2019animation, dipl\\u00f4me pr\\u00e9vu :\",\n    newText: \"Master of Arts in Animation, Attend graduation : Dec 2025\"\n  }\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the job-title/date lines for the three most recent positions and\n# the in-progress Master's degree line in the \"Exp\u00e9rience professionnelle\"\n# and \"Education\" sections of the r\u00e9sum\u00e9.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"Concepteur d\u2019animation (janvier 2021 \u00e0 ce jour)\"\n        New = \"Animation Spark : Concepteur d\u2019animations (Jan 2021 - Pr\u00e9sent)\"\n    },\n    @{\n        Old = \"Concepteur d\u2019animation (janvier 2018 \u00e0 d\u00e9cembre 2020)\"\n        New = \"Pixel Studio : Concepteur d\u2019animations (juin 2018 - D\u00e9c 2020)\"\n    },\n    @{\n        Old = \"Concepteur d\u2019animation junior (septembre 2016 \u00e0 mai 2018)\"\n        New = \"Animation flash : Concepteur d\u2019animation junior (sep 2016 - mai 2018)\"\n    },\n    @{\n        Old = \"Master en arts d\u2019animation, dipl\u00f4me pr\u00e9vu :\"\n        New = \"Master of Arts in Animation, Attend graduation : Dec 2025\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $found = $find.Execute($r.Old, $false, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)\n    if (-not $found) {\n        throw \"Text not found: $($r.Old)\"\n    }\n}\n"}
